$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data of row 2 and row 4 (observation records) for the
# columns that actually differ between them: A, B, E, F, G, H, Q, R, Z, AB.
# Row 3 and all other columns remain untouched.

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "Z", "AB")

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr4 = "$col" + "4"

    $val2 = $ws.Range($addr2).Value()
    $val4 = $ws.Range($addr4).Value()

    $ws.Range($addr2).Value = $val4
    $ws.Range($addr4).Value = $val2
}
